$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header of column B from seconds to minutes of sunshine
$ws.Range("B1").Value = "Solskensminuter_avg"

# Convert each value in column B (row 2 to 20) from seconds to minutes (rounded)
for ($r = 2; $r -le 20; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $v = $cell.Value()
    if ($v -ne $null) {
        $cell.Value = [Math]::Round($v / 60)
    }
}

# Update the selected cell/active cell on the sheet
$ws.Range("E5").Select()
